$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.343.13"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.937.18"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.73"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7248"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +2.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.96"
$ws.Range("E9").Value = "  +5.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07279"
$ws.Range("E10").Value = "  +6.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8097"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08105"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.934.81"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.488"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.03"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.17"
$ws.Range("E16").Value = "  +4.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.332.08"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008304"
$ws.Range("E18").Value = "  +6.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.93"
$ws.Range("E19").Value = "  -2.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.835"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.187.81"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("E22").Value = "  -0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.976"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.781"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.56"
$ws.Range("E26").Value = "  +4.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.353"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.37"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.357"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.540"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.447"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.217"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05261"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.272"
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7513"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.767"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("E38").Value = "  +3.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.807"
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "79.40"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.448"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4560"
$ws.Range("E42").Value = "  +3.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.038"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8465"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.94"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.468"
$ws.Range("E48").Value = "  +3.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.75"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4221"
$ws.Range("E50").Value = "  +4.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.502"
$ws.Range("E51").Value = "  +1.93%  "
